$p = $ppt.ActivePresentation

$oldDate = "5/29/2019"
$newDate = "7/28/2019"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master date placeholder
$master = $p.SlideMaster
Update-DateShapes $master.Shapes

# Every slide layout's date placeholder
$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $lyt = $layouts.Item($L)
    Update-DateShapes $lyt.Shapes
}

# Notes master date placeholder
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes
